$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 210, pushing existing rows 210-234 down to 212-236
$ws.Rows.Item(210).Resize(2).EntireRow.Insert()

# Fill in the new row 210 (Primera, 2000 kg, $/caja 10 kilos, Arica y Parinacota)
$ws.Cells.Item(210, 1).Value = 11
$ws.Cells.Item(210, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value = "Bíobío"
$ws.Cells.Item(210, 4).Value = 44449
$ws.Cells.Item(210, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(210, 5).Value = 8
$ws.Cells.Item(210, 6).Value = 100112020
$ws.Cells.Item(210, 7).Value = "Tomate"
$ws.Cells.Item(210, 8).Value = "Larga vida"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 2000
$ws.Cells.Item(210, 11).Value = 9000
$ws.Cells.Item(210, 12).Value = 10000
$ws.Cells.Item(210, 13).Value = 9500
$ws.Cells.Item(210, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(210, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(210, 16).Value = 950
$ws.Cells.Item(210, 17).Value = 10
$ws.Cells.Item(210, 18).Value = "Hortaliza"

# Fill in the new row 211 (Segunda, 1000 kg, $/caja 10 kilos, Arica y Parinacota)
$ws.Cells.Item(211, 1).Value = 11
$ws.Cells.Item(211, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(211, 3).Value = "Bíobío"
$ws.Cells.Item(211, 4).Value = 44449
$ws.Cells.Item(211, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(211, 5).Value = 8
$ws.Cells.Item(211, 6).Value = 100112020
$ws.Cells.Item(211, 7).Value = "Tomate"
$ws.Cells.Item(211, 8).Value = "Larga vida"
$ws.Cells.Item(211, 9).Value = "Segunda"
$ws.Cells.Item(211, 10).Value = 1000
$ws.Cells.Item(211, 11).Value = 8000
$ws.Cells.Item(211, 12).Value = 8000
$ws.Cells.Item(211, 13).Value = 8000
$ws.Cells.Item(211, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(211, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(211, 16).Value = 800
$ws.Cells.Item(211, 17).Value = 10
$ws.Cells.Item(211, 18).Value = "Hortaliza"
